$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: merged, centered header cell "temp" ---
$ws.Range("A1").Value = "temp"
$ws.Range("A1:B1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:B1").VerticalAlignment = -4108     # xlCenter
$ws.Range("A1:B1").Merge()

# --- Row 2 ---
$ws.Range("A2").Value = "dsr_pkt"
$ws.Range("B2").Value = "dsr基本数据包"

# --- Row 3 ---
$ws.Range("A3").Value = "rreq_tbl"
$ws.Range("B3").Value = "用来统计正在运行的路由发现，包括两种：1.经由本节点的 2.本节点发起的"

# --- Row 4 ---
$ws.Range("A4").Value = "LC"
$ws.Range("B4").Value = "Blacklist：没有函数会向其中添加内容，特殊需要时满足特殊功能（对所有被调用的函数而言）"

# --- Column B width (~74.625 characters) ---
$ws.Columns("B:B").ColumnWidth = 73.857142857

# --- Outline grouping: rows 1-4 nested 3 levels deep, column B grouped once ---
$ws.Rows("1:4").Group()
$ws.Rows("1:4").Group()
$ws.Rows("1:4").Group()
$ws.Columns("B:B").Group()

# --- Selection after edit ---
$ws.Range("B5").Select() | Out-Null
